# Elimina antiguos EC y agrega nuevos y modifica Antigua BD
# For this workbook: update "Periodo Mora" value from 2508 to 2509 for every
# worker row, and center-align that column's values (matching the rest of
# the table's centered columns).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Periodo Mora" column (E) for the three worker rows (16, 17, 18).
$periodoRange = $ws.Range("E16:E18")

# Update the period value: 2508 -> 2509 (kept as text, matching the
# existing "@" text number format already applied to this column).
$periodoRange.Value = "2509"

# Center the column horizontally, matching the rest of the data table.
$periodoRange.HorizontalAlignment = -4108  # xlCenter
